$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2:D10 currently hold "12.06.11" (shared-string text) -> change to "12.07.11"
# D11:D18 currently hold "12.07.11" (shared-string text) -> change to "12.08.11"
#
# These look like dates, so a plain `.Value = "12.07.11"` gets auto-parsed into
# a date serial by the COM layer. To keep them as literal text (matching the
# original cell typing / default style exactly, with no NumberFormat change
# recorded against the destination cells), stage the text in a same-shaped,
# text-formatted helper range off to the side, copy it, and paste *values only*
# into the destination. The helper range is cleared (not just emptied) afterward
# so it leaves no trace in the sheet's used range/dimension.

$xlPasteValues = -4163

# --- D2:D10 -> "12.07.11" ---
$helper1 = $ws.Range("Z1:Z9")
$helper1.NumberFormat = "@"
$helper1.Value = "12.07.11"
$helper1.Copy()
$ws.Range("D2:D10").PasteSpecial($xlPasteValues)
$helper1.Clear()

# --- D11:D18 -> "12.08.11" ---
$helper2 = $ws.Range("Z1:Z8")
$helper2.NumberFormat = "@"
$helper2.Value = "12.08.11"
$helper2.Copy()
$ws.Range("D11:D18").PasteSpecial($xlPasteValues)
$helper2.Clear()

# Update the sheet's active selection to D3:D10
$ws.Range("D3:D10").Select()
